$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 792.6667
$ws.Range("I18").Value = 798
$ws.Range("J18").Value = 790
$ws.Range("K18").Value = 798
$ws.Range("L18").Value = 790
$ws.Range("M18").Value = -514
$ws.Range("N18").Value = -1358
$ws.Range("H40").Value = 1405.7
$ws.Range("I40").Value = 1169.75
$ws.Range("J40").Value = 2349.5
$ws.Range("K40").Value = 1169.75
$ws.Range("L40").Value = 2349.5
$ws.Range("M40").Value = -994.75
$ws.Range("N40").Value = -2699.5
$ws.Range("H43").Value = 5789.25
$ws.Range("I43").Value = 19500
$ws.Range("J43").Value = 1219
$ws.Range("K43").Value = 19500
$ws.Range("L43").Value = 1219
$ws.Range("M43").Value = -19431
$ws.Range("N43").Value = -1357
$ws.Range("H51").Value = 2653.6667
$ws.Range("I51").Value = 2376.2856
$ws.Range("J51").Value = 3624.5
$ws.Range("K51").Value = 2376.2856
$ws.Range("L51").Value = 3624.5
$ws.Range("M51").Value = -1892.2856
$ws.Range("N51").Value = -4592.5
$ws.Range("H61").Value = 1375
$ws.Range("I61").Value = 1375
$ws.Range("K61").Value = 4125
$ws.Range("M61").Value = -3953
$ws.Range("H74").Value = 153785.28
$ws.Range("I74").Value = 254499.25
$ws.Range("J74").Value = 19500
$ws.Range("K74").Value = 254499.25
$ws.Range("L74").Value = 19500
$ws.Range("M74").Value = -253563.25
$ws.Range("N74").Value = -21372
$ws.Range("H77").Value = 153785.28
$ws.Range("I77").Value = 254499.25
$ws.Range("J77").Value = 19500
$ws.Range("K77").Value = 1272496.25
$ws.Range("L77").Value = 97500
$ws.Range("M77").Value = -1267816.25
$ws.Range("N77").Value = -106860
$ws.Range("H113").Value = 14322.182
$ws.Range("I113").Value = 12340.833
$ws.Range("K113").Value = 12340.833
$ws.Range("M113").Value = -9086.833000000001
$ws.Range("H116").Value = 2398.6667
$ws.Range("I116").Value = 2103.4
$ws.Range("J116").Value = 3875
$ws.Range("K116").Value = 2103.4
$ws.Range("L116").Value = 3875
$ws.Range("M116").Value = 1338.6
$ws.Range("N116").Value = -10759
$ws.Range("H138").Value = 6252537.5
$ws.Range("I138").Value = 1157.625
$ws.Range("J138").Value = 8931700
$ws.Range("K138").Value = 3472.875
$ws.Range("L138").Value = 26795100
$ws.Range("M138").Value = 1667.125
$ws.Range("N138").Value = -26805380

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 97499.5
$ws.Range("I45").Value = 97499.5
$ws.Range("K45").Value = 97499.5
$ws.Range("M45").Value = -97122.5
$ws.Range("H61").Value = 3666.6558
$ws.Range("I61").Value = 2464.8948
$ws.Range("K61").Value = 2464.8948
$ws.Range("M61").Value = -2252.8948
$ws.Range("H132").Value = 3224.3416
$ws.Range("I132").Value = 3060.162
$ws.Range("K132").Value = 9180.485999999999
$ws.Range("M132").Value = -6650.485999999999
$ws.Range("H136").Value = 3666.6558
$ws.Range("I136").Value = 2464.8948
$ws.Range("K136").Value = 7394.6844
$ws.Range("M136").Value = -4844.6844

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 1231.3334
$ws.Range("J12").Value = 1249.5
$ws.Range("L12").Value = 1249.5
$ws.Range("N12").Value = -1589.5
$ws.Range("H35").Value = 359.25
$ws.Range("I35").Value = 359.25
$ws.Range("K35").Value = 359.25
$ws.Range("M35").Value = -65.25
$ws.Range("H56").Value = 666.6667
$ws.Range("J56").Value = 500
$ws.Range("L56").Value = 500
$ws.Range("N56").Value = -2190
$ws.Range("H117").Value = 25000
$ws.Range("J117").Value = 25000
$ws.Range("L117").Value = 25000
$ws.Range("N117").Value = -34178
$ws.Range("H122").Value = 1733.44
$ws.Range("I122").Value = 1522.5555
$ws.Range("J122").Value = 2275.7144
$ws.Range("K122").Value = 4567.666499999999
$ws.Range("L122").Value = 6827.1432
$ws.Range("M122").Value = -2117.666499999999
$ws.Range("N122").Value = -11727.1432
$ws.Range("H132").Value = 3080.0588
$ws.Range("I132").Value = 3080.0588
$ws.Range("K132").Value = 9240.1764
$ws.Range("M132").Value = -6710.1764
$ws.Range("H134").Value = 13326.853
$ws.Range("I134").Value = 5770.467
$ws.Range("K134").Value = 17311.401
$ws.Range("M134").Value = -14776.401

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 191.22223
$ws.Range("J33").Value = 318
$ws.Range("L33").Value = 1908
$ws.Range("N33").Value = -2474
$ws.Range("H122").Value = 1746.0834
$ws.Range("I122").Value = 1872.6666
$ws.Range("J122").Value = 1366.3334
$ws.Range("K122").Value = 16853.9994
$ws.Range("L122").Value = 12297.0006
$ws.Range("M122").Value = -14403.9994
$ws.Range("N122").Value = -17197.0006

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 32000
$ws.Range("I46").Value = 8000
$ws.Range("J46").Value = 40000
$ws.Range("K46").Value = 8000
$ws.Range("L46").Value = 40000
$ws.Range("M46").Value = -7844
$ws.Range("N46").Value = -40312
$ws.Range("H80").Value = 4000
$ws.Range("J80").Value = 5500
$ws.Range("L80").Value = 5500
$ws.Range("N80").Value = -7496
$ws.Range("H83").Value = 4000
$ws.Range("J83").Value = 5500
$ws.Range("L83").Value = 27500
$ws.Range("N83").Value = -37484
$ws.Range("H97").Value = 759.41174
$ws.Range("I97").Value = 880.7143
$ws.Range("K97").Value = 880.7143
$ws.Range("M97").Value = -384.7143
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()
$ws.Range("H113").Value = 2766.2666
$ws.Range("I113").Value = 2813.8572
$ws.Range("K113").Value = 2813.8572
$ws.Range("M113").Value = -643.8571999999999
$ws.Range("H122").Value = 6996.6665
$ws.Range("I122").Value = 6000
$ws.Range("K122").Value = 18000
$ws.Range("M122").Value = -15550

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3319.111
$ws.Range("I22").Value = 3425
$ws.Range("J22").Value = 3234.4
$ws.Range("K22").Value = 3425
$ws.Range("L22").Value = 3234.4
$ws.Range("M22").Value = -3130
$ws.Range("N22").Value = -3824.4
$ws.Range("H27").Value = 3319.111
$ws.Range("I27").Value = 3425
$ws.Range("J27").Value = 3234.4
$ws.Range("K27").Value = 3425
$ws.Range("L27").Value = 3234.4
$ws.Range("M27").Value = -3318
$ws.Range("N27").Value = -3448.4
$ws.Range("H61").Value = 2351.7144
$ws.Range("I61").Value = 2453.76
$ws.Range("J61").Value = 1501.3334
$ws.Range("K61").Value = 2453.76
$ws.Range("L61").Value = 1501.3334
$ws.Range("M61").Value = -2251.76
$ws.Range("N61").Value = -1905.3334
$ws.Range("H113").Value = 2351.7144
$ws.Range("I113").Value = 2453.76
$ws.Range("J113").Value = 1501.3334
$ws.Range("K113").Value = 2453.76
$ws.Range("L113").Value = 1501.3334
$ws.Range("M113").Value = -283.7600000000002
$ws.Range("N113").Value = -5841.3334
$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").ClearContents()
$ws.Range("H122").Value = 441702.2
$ws.Range("I122").Value = 594432.7
$ws.Range("K122").Value = 1783298.1
$ws.Range("M122").Value = -1780848.1
$ws.Range("H132").Value = 4215.154
$ws.Range("I132").Value = 3816.4167
$ws.Range("K132").Value = 11449.2501
$ws.Range("M132").Value = -8919.250100000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H38").Value = 10056
$ws.Range("I38").Value = 10056
$ws.Range("K38").Value = 10056
$ws.Range("M38").Value = -9583
$ws.Range("H100").Value = 667.44446
$ws.Range("I100").Value = 647.8461
$ws.Range("K100").Value = 1295.6922
$ws.Range("M100").Value = -754.6922
$ws.Range("H113").Value = 1580.76
$ws.Range("I113").Value = 1548.2
$ws.Range("J113").Value = 1629.6
$ws.Range("K113").Value = 4644.6
$ws.Range("L113").Value = 4888.799999999999
$ws.Range("M113").Value = -2474.6
$ws.Range("N113").Value = -9228.799999999999
$ws.Range("H122").Value = 1762.3055
$ws.Range("I122").Value = 1647.7084
$ws.Range("J122").Value = 1991.5
$ws.Range("K122").Value = 4943.1252
$ws.Range("L122").Value = 5974.5
$ws.Range("M122").Value = -2493.1252
$ws.Range("N122").Value = -10874.5
$ws.Range("H135").Value = 74228.5
$ws.Range("J135").Value = 74228.5
$ws.Range("L135").Value = 74228.5
$ws.Range("N135").Value = -84368.5

